$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "5.6 [4.4-7.2] (n=102)" "5.8 [4.5-7.2] (n=101)"
Replace-Text "5.3 [5.0-7.4] (n=23)" "5.3 [4.9-7.4] (n=23)"
Replace-Text "0.06017575" "0.005951132"
Replace-Text "0.8062184" "0.93850933"
Replace-Text "7.1 [5.8-9.0] (n=102)" "7.0 [5.7-8.9] (n=100)"
Replace-Text "5.9 [5.2-7.7] (n=23)" "6.6 [5.4-9.2] (n=23)"
Replace-Text "2.14761805" "0.151478024"
Replace-Text "0.1427913" "0.69712689"
Replace-Text "3.2 [2.0-5.5] (n=107)" "3.4 [2.0-5.5] (n=106)"
Replace-Text "3.0 [1.0-4.0] (n=23)" "2.1 [0.2-3.9] (n=22)"
Replace-Text "1.50596582" "3.082514692"
Replace-Text "0.2197557" "0.07913808"
Replace-Text "26.3 [19.9-34.0] (n=92)" "26.5 [19.9-34.0] (n=91)"
Replace-Text "21.7 [17.8-31.1] (n=16)" "26.3 [19.6-31.8] (n=14)"
Replace-Text "0.54036199" "0.015017328"
Replace-Text "0.4622827" "0.90246725"
Replace-Text "5.9 [4.5-7.6] (n=104)" "5.9 [4.5-7.6] (n=103)"
Replace-Text "5.2 [4.4-7.5] (n=23)" "5.1 [4.3-7.5] (n=23)"
Replace-Text "0.25084641" "0.439756135"
Replace-Text "0.6164797" "0.50724018"
Replace-Text "7.8 [6.3-9.2] (n=102)" "7.8 [6.2-9.2] (n=101)"
Replace-Text "7.0 [6.0-8.8] (n=23)" "6.7 [6.1-8.8] (n=23)"
Replace-Text "0.30375650" "0.320037623"
Replace-Text "0.5815370" "0.57158504"
